$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.465.48'
$ws.Range('E2').Value = '  -2.05%  '

$ws.Range('D3').Value = '3.385.52'
$ws.Range('E3').Value = '  -1.67%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '567.44'
$ws.Range('E5').Value = '  -1.01%  '

$ws.Range('D6').Value = '159.57'
$ws.Range('E6').Value = '  +0.14%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = '3.386.45'
$ws.Range('E8').Value = '  -1.69%  '

$ws.Range('D9').Value = '0.547'
$ws.Range('E9').Value = '  -4.87%  '

$ws.Range('E10').Value = '  +0.87%  '

$ws.Range('D11').Value = '0.118'
$ws.Range('E11').Value = '  -2.97%  '

$ws.Range('D12').Value = '0.419'
$ws.Range('E12').Value = '  -5.13%  '

$ws.Range('D13').Value = '3.972.47'
$ws.Range('E13').Value = '  -1.59%  '

$ws.Range('E14').Value = '  +0.78%  '

$ws.Range('D15').Value = '26.73'
$ws.Range('E15').Value = '  -3.71%  '

$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').Value = '  -2.98%  '

$ws.Range('D17').Value = '63.570.99'
$ws.Range('E17').Value = '  -1.97%  '

$ws.Range('D18').Value = '3.303.46'
$ws.Range('E18').Value = '  -4.04%  '

$ws.Range('D19').Value = '6.06'
$ws.Range('E19').Value = '  -3.31%  '

$ws.Range('D20').Value = '13.47'
$ws.Range('E20').Value = '  -2.90%  '

$ws.Range('D21').Value = '373.34'
$ws.Range('E21').Value = '  -1.57%  '

$ws.Range('D22').Value = '7.70'
$ws.Range('E22').Value = '  -3.27%  '

$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.28%  '

$ws.Range('D24').Value = '70.46'
$ws.Range('E24').Value = '  -2.17%  '

$ws.Range('D25').Value = '0.510'
$ws.Range('E25').Value = '  -5.96%  '

$ws.Range('D26').Value = '0.0000113'
$ws.Range('E26').Value = '  -5.16%  '

$ws.Range('D27').Value = '9.38'
$ws.Range('E27').Value = '  -4.50%  '

$ws.Range('E28').Value = '  -0.03%  '

$ws.Range('E29').Value = '  +0.30%  '

$ws.Range('D30').Value = '5.96'
$ws.Range('E30').Value = '  -2.30%  '

$ws.Range('E31').Value = '  -7.12%  '

$ws.Range('D32').Value = '1.99'
$ws.Range('E32').Value = '  -0.84%  '

$ws.Range('D33').Value = '22.70'
$ws.Range('E33').Value = '  -2.07%  '

$ws.Range('D34').Value = '7.00'
$ws.Range('E34').Value = '  -1.08%  '

$ws.Range('D35').Value = '1.48'
$ws.Range('E35').Value = '  -5.65%  '

$ws.Range('D36').Value = '159.19'
$ws.Range('E36').Value = '  -1.24%  '

$ws.Range('D37').Value = '0.852'
$ws.Range('E37').Value = '  +9.41%  '

$ws.Range('D38').Value = '1.79'
$ws.Range('E38').Value = '  -5.10%  '

$ws.Range('D39').Value = '0.0716'
$ws.Range('E39').Value = '  -3.89%  '

$ws.Range('D40').Value = '2.747.36'
$ws.Range('E40').Value = '  -5.46%  '

$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '42.54'
$ws.Range('E41').Value = '  -0.88%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '25.42'
$ws.Range('E42').Value = '  -3.12%  '

$ws.Range('D43').Value = '26.01'
$ws.Range('E43').Value = '  -0.10%  '

$ws.Range('D44').Value = '6.35'
$ws.Range('E44').Value = '  -4.73%  '

$ws.Range('D45').Value = '4.33'
$ws.Range('E45').Value = '  -4.48%  '

$ws.Range('D46').Value = '0.0304'
$ws.Range('E46').Value = '  -2.56%  '

$ws.Range('D47').Value = '2.38'
$ws.Range('E47').Value = '  +4.60%  '

$ws.Range('D48').Value = '324.08'
$ws.Range('E48').Value = '  +2.22%  '

$ws.Range('E49').Value = '  -5.27%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.102'
$ws.Range('E50').Value = '  -2.77%  '

$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '6.25'
$ws.Range('E51').Value = '  -3.67%  '
